{"js": "/*\n * Update the date line and all 100 math-problem cells in the single\n * table. Each cell/paragraph keeps its existing formatting because we\n * call insertText(\"Replace\") on the *existing* paragraph range instead\n * of clearing/rebuilding the body (which would drop the run's rPr and\n * the paragraph's pPr).\n */\nconst newTitle = \"2025-05-02 Friday\";\n\n// New values for the 20x5 grid of practice problems, in row-major order\n// (matches the table's row/column layout: newGrid[row][col]).\nconst newGrid = [\n  [\"75-74=\", \"49-21=\", \"11+77=\", \"81-62=\", \"81-11=\"],\n  [\"97-7=\", \"77-17=\", \"69-12=\", \"71-15=\", \"60+35=\"],\n  [\"76-9=\", \"91-29=\", \"66-40=\", \"38+21=\", \"52+45=\"],\n  [\"71+26=\", \"70-33=\", \"85-58=\", \"42-36=\", \"31-24=\"],\n  [\"92-38=\", \"25+6=\", \"84-60=\", \"97-37=\", \"61-16=\"],\n  [\"74-64=\", \"19+64=\", \"22+42=\", \"47-4=\", \"25+17=\"],\n  [\"95-65=\", \"14+27=\", \"7+46=\", \"87-15=\", \"34-13=\"],\n  [\"87+5=\", \"57-15=\", \"87-76=\", \"12+22=\", \"45-23=\"],\n  [\"64+30=\", \"33-17=\", \"45-17=\", \"8+7=\", \"2+85=\"],\n  [\"22-14=\", \"86+11=\", \"21+16=\", \"92-68=\", \"11+45=\"],\n  [\"36+48=\", \"74-55=\", \"39+14=\", \"38+17=\", \"45+40=\"],\n  [\"32+26=\", \"23+61=\", \"76-37=\", \"71-52=\", \"46+43=\"],\n  [\"32+59=\", \"21-1=\", \"1+61=\", \"15-3=\", \"69-58=\"],\n  [\"6+20=\", \"39+38=\", \"10-1=\", \"13+83=\", \"78+0=\"],\n  [\"21-15=\", \"87-58=\", \"72+6=\", \"23+75=\", \"13+23=\"],\n  [\"29+68=\", \"22-18=\", \"4-4=\", \"49-17=\", \"38-31=\"],\n  [\"22+56=\", \"61-39=\", \"22-9=\", \"34+40=\", \"99-58=\"],\n  [\"79+5=\", \"66+4=\", \"59+29=\", \"88-74=\", \"68-29=\"],\n  [\"39+40=\", \"68+15=\", \"27-3=\", \"34-15=\", \"46-2=\"],\n  [\"24+39=\", \"15+4=\", \"73-31=\", \"71+6=\", \"62-24=\"]\n];\n\nconst body = context.document.body;\n\n// --- Update the title/date paragraph (first paragraph in the body) ---\nconst titlePara = body.paragraphs.getFirst();\ntitlePara.insertText(newTitle, \"Replace\");\n\n// --- Update every cell of the first (only) table ---\nconst table = body.tables.getFirst();\nfor (let r = 0; r < newGrid.length; r++) {\n  const rowValues = newGrid[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cellPara = table.getCell(r, c).body.paragraphs.getFirst();\n    cellPara.insertText(rowValues[c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Update the date/title line (first paragraph) ---\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleRange.End = $titleRange.End - 1\n$titleRange.Text = \"2025-05-02 Friday\"\n\n# --- Update every cell of the single practice-problem table ---\n$newValues = @(\n    @(\"75-74=\", \"49-21=\", \"11+77=\", \"81-62=\", \"81-11=\"),\n    @(\"97-7=\", \"77-17=\", \"69-12=\", \"71-15=\", \"60+35=\"),\n    @(\"76-9=\", \"91-29=\", \"66-40=\", \"38+21=\", \"52+45=\"),\n    @(\"71+26=\", \"70-33=\", \"85-58=\", \"42-36=\", \"31-24=\"),\n    @(\"92-38=\", \"25+6=\", \"84-60=\", \"97-37=\", \"61-16=\"),\n    @(\"74-64=\", \"19+64=\", \"22+42=\", \"47-4=\", \"25+17=\"),\n    @(\"95-65=\", \"14+27=\", \"7+46=\", \"87-15=\", \"34-13=\"),\n    @(\"87+5=\", \"57-15=\", \"87-76=\", \"12+22=\", \"45-23=\"),\n    @(\"64+30=\", \"33-17=\", \"45-17=\", \"8+7=\", \"2+85=\"),\n    @(\"22-14=\", \"86+11=\", \"21+16=\", \"92-68=\", \"11+45=\"),\n    @(\"36+48=\", \"74-55=\", \"39+14=\", \"38+17=\", \"45+40=\"),\n    @(\"32+26=\", \"23+61=\", \"76-37=\", \"71-52=\", \"46+43=\"),\n    @(\"32+59=\", \"21-1=\", \"1+61=\", \"15-3=\", \"69-58=\"),\n    @(\"6+20=\", \"39+38=\", \"10-1=\", \"13+83=\", \"78+0=\"),\n    @(\"21-15=\", \"87-58=\", \"72+6=\", \"23+75=\", \"13+23=\"),\n    @(\"29+68=\", \"22-18=\", \"4-4=\", \"49-17=\", \"38-31=\"),\n    @(\"22+56=\", \"61-39=\", \"22-9=\", \"34+40=\", \"99-58=\"),\n    @(\"79+5=\", \"66+4=\", \"59+29=\", \"88-74=\", \"68-29=\"),\n    @(\"39+40=\", \"68+15=\", \"27-3=\", \"34-15=\", \"46-2=\"),\n    @(\"24+39=\", \"15+4=\", \"73-31=\", \"71+6=\", \"62-24=\")\n)\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $cellRange = $table.Cell($r, $c).Range\n        $cellRange.End = $cellRange.End - 2\n        $cellRange.Text = $rowValues[$c - 1]\n    }\n}\n"}
